$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark that Word inserts to mark the last edit
#    position; later saves from this document no longer carry it.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# 2. Merge the ">>> your stuff after this line >>>" paragraph back into a
#    single run (the original had been split up by spell/grammar-check
#    proofing marks around "your").
$d.Content.Find.Execute(">>>  your stuff after this line >>>", $false, $false, $false, $false, $false, $true, 1, $false, ">>>  your stuff after this line >>>", 2)

# 3. Replace "Ben changing things up!" with "TaorongTang changing things up!"
$d.Content.Find.Execute("Ben changing things up!", $false, $false, $false, $false, $false, $true, 1, $false, "TaorongTang changing things up!", 2)

# Split "TaorongTang" into its own run (distinct from " changing things up!")
$f = $d.Content
$f.Find.Execute("TaorongTang")
$f.Bold = $true
$f.Bold = $false

Write-Output "done"
